$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New date-as-text values (DD/MM/YYYY), replacing the old date-serial values in D2:D13
$dates = @(
    "11/03/2013",
    "14/05/2013",
    "17/07/2013",
    "19/09/2013",
    "22/11/2013",
    "25/01/2014",
    "30/03/2014",
    "02/06/2014",
    "05/08/2014",
    "08/10/2014",
    "11/12/2014",
    "13/02/2015"
)

# Delete the now-empty extra rows 14:23 first (so the used range/dimension shrinks)
$ws.Range("A14:E23").EntireRow.Delete() | Out-Null

# Change the number format for column D (style used by header + data cells) from date to text
$ws.Range("D1:D13").NumberFormat = "@"

# Write the text date values into D2:D13
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $dates[$i]
}

# Update the visible selection to match the target state
$ws.Range("D14:D20").Select() | Out-Null
